# Training Dashboard update: progress refreshed as of 04-Nov-2025.
# - Rows 3-20: "LAST UPDATE" (I) moves from 03-Nov-2025 to 04-Nov-2025 and
#   "PERIOD TO EXPIRE" (H) decreases accordingly (rows 12 & 13 additionally
#   got a brand-new training/expiry date pair).
# - A brand-new training ("Diagnosis Of Beckoff Module...") is inserted as
#   new row 21, pushing every following row down by one.
# - All the pushed-down rows (now 22-36) get their serial number (A),
#   period-to-expire (H) and last-update (I) refreshed too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Writing a date-looking string via .Value lets Excel auto-convert it
    # into a real date serial. A leading apostrophe forces it to stay text,
    # exactly like typing it in the Excel UI would.
    $range.Value = "'" + $text
}

# ---------------------------------------------------------------------
# 1) Refresh the existing rows 3-20 (H = days left, I = last update date)
# ---------------------------------------------------------------------
$updates_3_20 = @(
    @{ Row = 3;  H = 703 },
    @{ Row = 4;  H = 704 },
    @{ Row = 5;  H = 712 },
    @{ Row = 6;  H = 702 },
    @{ Row = 7;  H = 712 },
    @{ Row = 8;  H = 368 },
    @{ Row = 9;  H = 704 },
    @{ Row = 10; H = 712 },
    @{ Row = 11; H = 703 },
    @{ Row = 12; H = 714; F = "20-Oct-2025"; G = "20-Oct-2027" },
    @{ Row = 13; H = 714; F = "20-Oct-2025"; G = "20-Oct-2027" },
    @{ Row = 14; H = 361 },
    @{ Row = 15; H = 362 },
    @{ Row = 16; H = 705 },
    @{ Row = 17; H = 425 },
    @{ Row = 18; H = 424 },
    @{ Row = 19; H = 423 },
    @{ Row = 20; H = 424 }
)

foreach ($u in $updates_3_20) {
    $r = $u.Row
    if ($u.ContainsKey("F")) {
        Set-TextValue $ws.Range("F$r") $u.F
        Set-TextValue $ws.Range("G$r") $u.G
    }
    $ws.Range("H$r").Value = $u.H
    Set-TextValue $ws.Range("I$r") "04-Nov-2025"
}

# ---------------------------------------------------------------------
# 2) Insert the new training row at position 21, cloning the formatting
#    of the row right above it (row 20) so borders/fill/alignment match.
# ---------------------------------------------------------------------
$ws.Rows.Item(21).Insert()
$ws.Range("A20:K20").Copy()
$ws.Range("A21:K21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "Diagnosis Of Beckoff Module And Troubleshooting Guide (Cargo Trainings)"
$ws.Range("C21").Value = "CARGO"
$ws.Range("D21").Value = "LSME-CRG-M-012"
$ws.Range("E21").Value = "EQUIPMENT MANUAL"
Set-TextValue $ws.Range("F21") "21-Oct-2025"
Set-TextValue $ws.Range("G21") "21-Oct-2027"
$ws.Range("H21").Value = 715
Set-TextValue $ws.Range("I21") "04-Nov-2025"
$ws.Range("J21").Value = "VALID"
$ws.Range("K21").Value = ""

# ---------------------------------------------------------------------
# 3) The rows that used to be 21-35 are now 22-36 (content/formatting was
#    already shifted down by the Insert). Refresh their serial number (A),
#    period-to-expire (H) and last-update (I).
# ---------------------------------------------------------------------
$updates_22_36 = @(
    @{ Row = 22; A = 20; H = 35 },
    @{ Row = 23; A = 21; H = -144 },
    @{ Row = 24; A = 22; H = -104 },
    @{ Row = 25; A = 23; H = 137 },
    @{ Row = 26; A = 24; H = 136 },
    @{ Row = 27; A = 25; H = 151 },
    @{ Row = 28; A = 26; H = 151 },
    @{ Row = 29; A = 27; H = 263 },
    @{ Row = 30; A = 28; H = 263 },
    @{ Row = 31; A = 29; H = 263 },
    @{ Row = 32; A = 30; H = 263 },
    @{ Row = 33; A = 31; H = 347 },
    @{ Row = 34; A = 32; H = 284 },
    @{ Row = 35; A = 33; H = 284 },
    @{ Row = 36; A = 34; H = 703 }
)

foreach ($u in $updates_22_36) {
    $r = $u.Row
    $ws.Range("A$r").Value = $u.A
    $ws.Range("H$r").Value = $u.H
    Set-TextValue $ws.Range("I$r") "04-Nov-2025"
}
